$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text helper: assign a value to a cell while guaranteeing it stays
# a Text-typed cell (matches the source inlineStr convention) even when the
# string looks numeric (e.g. "167.00"), and leaves the cell style untouched.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "35.398.28"
$ws.Range("E2").Value = "  +1.73%  "
Set-TextValue $ws.Range("D3") "1.888.86"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue $ws.Range("D5") "246.41"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  +1.27%  "
$ws.Range("E7").Value = "  +0.01%  "
Set-TextValue $ws.Range("D8") "43.23"
$ws.Range("E8").Value = "  +3.59%  "
$ws.Range("E9").Value = "  +3.06%  "
Set-TextValue $ws.Range("D10") "54.84"
$ws.Range("E10").Value = "  +7.56%  "
Set-TextValue $ws.Range("D11") "0.0745"
$ws.Range("E11").Value = "  +1.30%  "
Set-TextValue $ws.Range("D12") "0.0988"
$ws.Range("E12").Value = "  +1.78%  "
Set-TextValue $ws.Range("D13") "13.97"
$ws.Range("E13").Value = "  +8.49%  "
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("E15").Value = "  +8.47%  "
$ws.Range("E16").Value = "  +2.85%  "
Set-TextValue $ws.Range("D17") "1.879.86"
$ws.Range("E17").Value = "  +0.25%  "
Set-TextValue $ws.Range("D18") "35.432.75"
$ws.Range("E18").Value = "  +1.78%  "
Set-TextValue $ws.Range("D19") "73.49"
$ws.Range("E19").Value = "  +1.16%  "
Set-TextValue $ws.Range("D20") "0.0₃0826"
$ws.Range("E20").Value = "  +1.02%  "
Set-TextValue $ws.Range("D21") "245.44"
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("E23").Value = "  +5.10%  "
Set-TextValue $ws.Range("D24") "2.67"
$ws.Range("E24").Value = "  +8.83%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -3.41%  "
Set-TextValue $ws.Range("D27") "167.00"
$ws.Range("E27").Value = "  +1.24%  "
Set-TextValue $ws.Range("D28") "8.59"
$ws.Range("E28").Value = "  +2.83%  "
Set-TextValue $ws.Range("D29") "18.33"
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("E30").Value = "  +0.90%  "
Set-TextValue $ws.Range("D31") "0.0596"
$ws.Range("E31").Value = "  +3.49%  "
Set-TextValue $ws.Range("D32") "4.29"
$ws.Range("E32").Value = "  +1.90%  "
$ws.Range("E33").Value = "  +25.94%  "
$ws.Range("E34").Value = "  +1.94%  "
$ws.Range("E35").Value = "  -0.04%  "
Set-TextValue $ws.Range("D36") "1.47"
$ws.Range("E36").Value = "  -12.70%  "
Set-TextValue $ws.Range("D37") "0.857"
$ws.Range("E37").Value = "  +3.67%  "
$ws.Range("E38").Value = "  -1.71%  "
Set-TextValue $ws.Range("D39") "0.0720"
$ws.Range("E39").Value = "  +8.75%  "
Set-TextValue $ws.Range("D40") "0.0223"
$ws.Range("E40").Value = "  +6.72%  "
Set-TextValue $ws.Range("D41") "98.54"
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("E42").Value = "  +1.72%  "
$ws.Range("E43").Value = "  -0.42%  "
Set-TextValue $ws.Range("D44") "13.78"
$ws.Range("E44").Value = "  +15.14%  "
Set-TextValue $ws.Range("D45") "1.327.40"
$ws.Range("E45").Value = "  +3.48%  "
Set-TextValue $ws.Range("D46") "2.40"
$ws.Range("E46").Value = "  +3.33%  "
Set-TextValue $ws.Range("D47") "0.0810"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("E50").Value = "  -1.69%  "

# Row 51 full replacement: RocketPoolETH -> MultiversX
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue $ws.Range("D51") "42.31"
$ws.Range("E51").Value = "  -1.38%  "
